$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Update the JOIN note for the DimPatient -> PatientSurvey rows (Height(cms), Height(In),
# Weight(Lbs), Weight(Kgs), BloodGroup, Tobacco, Alcohol, Exercise, Diet, Ethinicity):
# the join column used for the "Spaghetti" lookup changed from HealthCardNbr to
# PatientNbr/SurveyNbr.
$oldNote = "JOIN using Patient.HealthCardNbr = PatientSurvey.HealthCardNbr"
$newNote = "JOIN using Patient.PatientNbr = PatientSurvey.SurveyNbr"

for ($row = 24; $row -le 33; $row++) {
    $cell = $ws.Range("E$row")
    if ($cell.Value2 -eq $oldNote) {
        $cell.Value = $newNote
    }
}

# Restore the view/selection state recorded for the sheet after the edit.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G29").Select()
